$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# Insert two new columns at P:Q, shifting everything from P onward to the right.
$ws.Range("P1:Q1").EntireColumn.Insert()

Write-Host ("Dimension: " + $ws.UsedRange.Address())

$ws.Range("P6:Q10").NumberFormat = "yyyy/mm/dd"
Write-Host ("P6 numfmt: " + $ws.Range("P6").NumberFormat)

$ws.Range("P5").Value = "DOI"
$ws.Range("Q5").Value = "DOE"
$ws.Range("O5").Value = "code"

$ws.Range("P6").Formula = "=NOW()-20000"
$ws.Range("Q6").Formula = "=NOW()+20000"
$ws.Range("P7").Formula = "=NOW()-20000"
$ws.Range("Q7").Formula = "=NOW()+20000"
$ws.Range("P7:P10").FillDown()
$ws.Range("Q7:Q10").FillDown()



